$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values (Default/Country/Sponsor specific sites sync; institute fields added) ---
$ws.Range("A2").Value = "Principal"
$ws.Range("B2").Value = "0102/0304"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = "DeLuca Jr., William F MD"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "DeLuca"
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = "William"
$ws.Range("J2").Value = "Apollo"
$ws.Range("K2").Value = "abc street"
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = "NY"
$ws.Range("N2").Value = "US"
$ws.Range("O2").Value = 889
$ws.Range("P2").Value = "Spain"
$ws.Range("Q2").Value = "#6789"

# --- Formatting: group of cells using Calibri (explicit) center/center + wrap ---
foreach ($addr in @("A2","C2","D2","E2","F2","G2")) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Calibri"
    $r.Font.Size = 11
    $r.Font.ThemeColor = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
    $r.WrapText = $true
}

# --- Formatting: group using Arial 10 (existing font) center/center, no wrap ---
foreach ($addr in @("H2","J2","K2")) {
    $r = $ws.Range($addr)
    $r.Font.Name = "Arial"
    $r.Font.Size = 10
    $r.Font.Color = 0
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}

# --- Formatting: group using default font, center/center, no wrap ---
foreach ($addr in @("B2","I2","L2","N2","O2","Q2")) {
    $r = $ws.Range($addr)
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4108
}

# --- Formatting: M2 uses default font, center/center, with wrap ---
$m2 = $ws.Range("M2")
$m2.HorizontalAlignment = -4108
$m2.VerticalAlignment = -4108
$m2.WrapText = $true

# --- Formatting: P2 uses Arial 11 color 333333, no special alignment ---
$p2 = $ws.Range("P2")
$p2.Font.Name = "Arial"
$p2.Font.Size = 11
$p2.Font.Color = 3355443

# --- Update selection to A2 ---
$ws.Range("A2").Select()
